$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: make a header row bold / centered / top-aligned / boxed, matching
# the look of the header rows already present in the workbook.
# ---------------------------------------------------------------------------
function Format-HeaderRow($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous
}

# ---------------------------------------------------------------------------
# Locate the two pre-existing sheets (their internal order never changes,
# only their position within the workbook tab-strip does).
# ---------------------------------------------------------------------------
$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")

# ===========================================================================
# 1) "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (header + values),
#    and drop the stray empty INNING_NUMBER cells.
# ===========================================================================
$battingWs.Range("D1").Value = "MATCH_CODE"

$lastRow = $battingWs.Cells.Item($battingWs.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $linkCell = $battingWs.Cells.Item($r, 4)
    $link = $linkCell.Value2
    if ($link -ne $null -and $link -ne "") {
        $code = [string]$link -replace '.*MatchCode=', ''
        $linkCell.Value = "'" + $code
    }

    $inningCell = $battingWs.Cells.Item($r, 2)
    $inning = $inningCell.Value2
    if ($inning -eq $null -or $inning -eq "") {
        $inningCell.ClearContents()
    }
}

# ===========================================================================
# 2) "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (header + values).
# ===========================================================================
$bowlingWs.Range("B1").Value = "MATCH_CODE"

$lastRow = $bowlingWs.Cells.Item($bowlingWs.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $linkCell = $bowlingWs.Cells.Item($r, 2)
    $link = $linkCell.Value2
    if ($link -ne $null -and $link -ne "") {
        $code = [string]$link -replace '.*MatchCode=', ''
        $linkCell.Value = "'" + $code
    }
}

# ===========================================================================
# 3) New "Player Info" sheet, inserted as the very first tab.
# ===========================================================================
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
Format-HeaderRow $playerInfo.Range("A1:D1")

$playerInfo.Range("A2").Value = "'3788"
$playerInfo.Range("B2").Value = "Ravichandran Ashwin"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# ===========================================================================
# 4) New "ODI Batting Extra" sheet, appended as the very last tab.
# ===========================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
Format-HeaderRow $extra.Range("A1:F1")

# Each data line: MATCH_CODE|BATTING_POSITION|NUM_4|NUM_6|PERCENT_RUNS_OF_TOTAL|MAN_OF_MATCH
# A numeric BATTING_POSITION is prefixed with '#'; blank fields are left untouched.
$extraData = @"
3786|||||NO
3791|#8|0|0|0.99%|NO
3795|#8|0|0|2.15%|NO
3808|#8|0|0||NO
3810|||||NO
3811|#9||||NO
3841|#10||||NO
3874|#8||||NO
3875|||||NO
3974|#9|0|1|4.21%|NO
3976|||||NO
3978|||||NO
4042|#10||||NO
4047|#9||||NO
4050|||||NO
4051|||||NO
4052|#8||||NO
4053|||||NO
4524|||||NO
4526|#8|1|1|8.71%|NO
"@

$lines = $extraData -split "`n"
$rowIdx = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split '\|', -1

    $extra.Range("A$rowIdx").Value = "'" + $parts[0]

    if ($parts[1] -ne "") {
        $extra.Range("B$rowIdx").Value = [double]($parts[1].Substring(1))
    }
    if ($parts[2] -ne "") {
        $extra.Range("C$rowIdx").Value = "'" + $parts[2]
    }
    if ($parts[3] -ne "") {
        $extra.Range("D$rowIdx").Value = "'" + $parts[3]
    }
    if ($parts[4] -ne "") {
        $extra.Range("E$rowIdx").Value = "'" + $parts[4]
    }
    $extra.Range("F$rowIdx").Value = $parts[5]

    $rowIdx++
}

Write-Output "Final sheet order:"
foreach ($s in $wb.Worksheets) { Write-Output $s.Name }
